# Updated confidence-interval / estimate values feeding the plot on the
# "For plotting" sheet (commit message: "changed plot axis limits").
# The underlying re-computed statistics changed slightly, which moves the
# bootstrap CI values referenced by the plot's axis limits.

$wb = $excel.ActiveWorkbook

$wsFull = $wb.Worksheets.Item("Full results")
$wsPlot = $wb.Worksheets.Item("For plotting")

# --- "Full results" sheet -------------------------------------------------

# Row 2 (education, COMPLETE MODEL)
$wsFull.Range("H2").Value = 0.563380826992887
$wsFull.Range("I2").Value = 0.180654977621478
$wsFull.Range("O2").Value = 0.436670145128977

# Row 3 (education, CONDITIONAL MODEL)
$wsFull.Range("F3").Value = 0.575939809378633
$wsFull.Range("G3").Value = 0.203514636764776

# Row 4 (education, NULL MODEL)
$wsFull.Range("C4").Value = 0.63758128463405
$wsFull.Range("D4").Value = 0.362509198923953
$wsFull.Range("E4").Value = 1.000090483558
$wsFull.Range("J4").Value = 0.362476400765031
$wsFull.Range("K4").Value = 0.203496223722589
$wsFull.Range("L4").Value = 0.0125578461080231
$wsFull.Range("M4").Value = 0.0741937443639465
$wsFull.Range("N4").Value = 0.216054069830612

# --- "For plotting" sheet --------------------------------------------------

# Row 2 (IOLIB)
$wsPlot.Range("C2").Value = 0.362476400765031
$wsPlot.Range("D2").Value = 0.330148319043034
$wsPlot.Range("E2").Value = 0.394804482487027

# Row 3 (IORAD)
$wsPlot.Range("C3").Value = 0.216054069830612
$wsPlot.Range("D3").Value = 0.188708379640886
$wsPlot.Range("E3").Value = 0.243399760020337

# Row 4 (Sibcorr)
$wsPlot.Range("C4").Value = 0.436670145128977
$wsPlot.Range("D4").Value = 0.405504453035125
$wsPlot.Range("E4").Value = 0.46783583722283
